# Apply updated vm_pu results ("case with 380 kV done") to Sheet1.
# Rows 2-25 correspond to buses 0-23; columns B-F and I-N hold the
# per-scenario voltage magnitude (p.u.) results that were recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ 2=1.02; 3=1.022613565826668; 4=1.032304562201916; 5=0.9926147277508489; 6=1.041002484908686; 9=1.032048314739603; 10=1.027798416421324; 11=1.03511017814932; 12=0.9955398523336033; 13=1.043783233570425; 14=1.013331343508064 }
    3 = @{ 2=1.02; 3=1.023396567042635; 4=1.032882560340264; 5=0.9936372048519304; 6=1.041717966528747; 9=1.03214341065381; 10=1.028220406776442; 11=1.035497538398203; 12=0.9963617723202692; 13=1.044309503773916; 14=1.013471714926972 }
    4 = @{ 2=1.02; 3=1.023903894802248; 4=1.03325699817294; 5=0.9942998659930995; 6=1.042181711634109; 9=1.032203790458801; 10=1.028493471037827; 11=1.035747919490475; 12=0.9968940712668345; 13=1.044650143282474; 14=1.013562520721699 }
    5 = @{ 2=1.02; 3=1.024117334797243; 4=1.033414513859839; 5=0.9945786998346017; 6=1.042376854771378; 9=1.032228897230713; 10=1.028608267894394; 11=1.035853114378514; 12=0.997117960005301; 13=1.044793372094303; 14=1.013600689279714 }
    6 = @{ 2=1.02; 3=1.024153181605111; 4=1.033440967361991; 5=0.9946255319796338; 6=1.042409630929501; 9=1.032233096509786; 10=1.028627542790286; 11=1.035870773190927; 12=0.9971555583673453; 13=1.044817422195631; 14=1.013607097571227 }
    7 = @{ 2=1.02; 3=1.023906746173885; 4=1.033259102506866; 5=0.9943035907982488; 6=1.0421843184233; 9=1.032204127025637; 10=1.028495004958211; 11=1.03574932536819; 12=0.9968970624462087; 13=1.044652057020943; 14=1.013563030756792 }
    8 = @{ 2=1.02; 3=1.022878044158395; 4=1.032499808252087; 5=0.9929600610674301; 6=1.041244122313313; 9=1.032080691277773; 10=1.027941027515717; 11=1.035241142761852; 12=0.995817528259106; 13=1.043961065726742; 14=1.013378787348611 }
    9 = @{ 2=1.02; 3=1.021070581174474; 4=1.03116524871004; 5=0.9906006454969559; 6=1.039593452524476; 9=1.031854378358587; 10=1.026964971241679; 11=1.034343679151538; 12=0.9939188001724441; 13=1.042744352385792; 14=1.013053964926507 }
    10 = @{ 2=1.02; 3=1.019869240967673; 4=1.030277955418011; 5=0.989033133672735; 6=1.038497222254922; 9=1.031697627330683; 10=1.026314429718382; 11=1.033744123536263; 12=0.9926553831429383; 13=1.041933922211437; 14=1.012837336154905 }
    11 = @{ 2=1.02; 3=1.019349933065729; 4=1.029894344929067; 5=0.988355674866747; 6=1.038023570292174; 9=1.031628367870288; 10=1.026032793586741; 11=1.033484231732386; 12=0.9921088820399291; 13=1.041583188310763; 14=1.01274352042715 }
    12 = @{ 2=1.02; 3=1.019157173146807; 4=1.029751946090618; 5=0.9881042295826724; 6=1.037847790820111; 9=1.031602434430107; 10=1.025928190667124; 11=1.033387655697548; 12=0.9919059725120875; 13=1.041452940144107; 14=1.012708671465653 }
    13 = @{ 2=1.02; 3=1.019198514702341; 4=1.029782486966316; 5=0.9881581567098651; 6=1.037885489019889; 9=1.031608006614418; 10=1.025950627906556; 11=1.033408373400817; 12=0.9919494934313052; 13=1.041480877430334; 14=1.01271614675184 }
    14 = @{ 2=1.02; 3=1.019333996722937; 4=1.029882572324939; 5=0.9883348863814464; 6=1.038009037111634; 9=1.031626228430603; 10=1.026024146879554; 11=1.033476249547978; 12=0.9920921077337197; 13=1.041572421333761; 14=1.012740639831422 }
    15 = @{ 2=1.02; 3=1.019417489501323; 4=1.029944250360185; 5=0.9884438009545853; 6=1.038085179901962; 9=1.031637428033007; 10=1.026069445615852; 11=1.033518064939908; 12=0.9921799884222134; 13=1.041628828586344; 14=1.012755730625585 }
    16 = @{ 2=1.02; 3=1.019903724432428; 4=1.030303427057747; 5=0.9890781214508737; 6=1.038528678722194; 9=1.031702194729407; 10=1.026333122207026; 11=1.033761365919115; 12=0.9926916645766087; 13=1.041957203375523; 14=1.012843562141503 }
    17 = @{ 2=1.02; 3=1.020208963788498; 4=1.0305288894063; 5=0.989476357848556; 6=1.038807149295824; 9=1.031742450811208; 10=1.026498534628113; 11=1.033913908191406; 12=0.9930127773699352; 13=1.042163235746173; 14=1.012898653115598 }
    18 = @{ 2=1.02; 3=1.020387089552913; 4=1.030660454894529; 5=0.9897087662937556; 6=1.038969675020917; 9=1.031765797799479; 10=1.026595021882016; 11=1.034002856378659; 12=0.9932001317071769; 13=1.042283428931943; 14=1.012930785361489 }
    19 = @{ 2=1.02; 3=1.020447840163313; 4=1.030705324951847; 5=0.9897880325774034; 6=1.039025108740085; 9=1.031773735818787; 10=1.026627922346335; 11=1.034033180764138; 12=0.9932640239640975; 13=1.042324414698219; 14=1.012941741376859 }
    20 = @{ 2=1.02; 3=1.020176205708132; 4=1.030504693507743; 5=0.9894336180360679; 6=1.038777261828776; 9=1.031738145538948; 10=1.026480786915036; 11=1.033897544645857; 12=0.9929783193494215; 13=1.042141128548788; 14=1.012892742515433 }
    21 = @{ 2=1.02; 3=1.019294096926008; 4=1.029853097148155; 5=0.9882828385668249; 6=1.037972650954643; 9=1.031620868279157; 10=1.026002497107258; 11=1.03345626282476; 12=0.9920501090198102; 13=1.041545463092571; 14=1.012733427270863 }
    22 = @{ 2=1.02; 3=1.018740257261202; 4=1.029443941181273; 5=0.9875604150241495; 6=1.037467663596203; 9=1.031545931858191; 10=1.025701831687648; 11=1.033178576978026; 12=0.9914670000341481; 13=1.041171119165601; 14=1.012633250302939 }
    23 = @{ 2=1.02; 3=1.019033783848894; 4=1.029660791725634; 5=0.9879432794643023; 6=1.037735280559257; 9=1.03158577053523; 10=1.02586121450481; 11=1.033325805235903; 12=0.991776070289318; 13=1.041369548808774; 14=1.012686356714251 }
    24 = @{ 2=1.02; 3=1.020191007412783; 4=1.030515626416597; 5=0.9894529299347244; 6=1.03879076638444; 9=1.031740091319377; 10=1.026488806327404; 11=1.033904938711777; 12=0.9929938892766442; 13=1.042151117784347; 14=1.012895413265579 }
    25 = @{ 2=1.02; 3=1.021537221105039; 4=1.031509847622851; 5=0.9912096547607049; 6=1.040019456677264; 9=1.031913924473337; 10=1.027217282375347; 11=1.034575920957438; 12=0.9944092447426414; 13=1.043058784185277; 14=1.013137955573838 }
}

foreach ($rowKey in $newValues.Keys) {
    $rowData = $newValues[$rowKey]
    foreach ($colKey in $rowData.Keys) {
        $ws.Cells.Item([int]$rowKey, [int]$colKey).Value = $rowData[$colKey]
    }
}
